$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = New-Object 'object[,]' 1,189
$values[0,0] = [double]"2.055726895378029e-07"
$values[0,1] = [double]"1.497821267548716e-05"
$values[0,2] = [double]"4.106057986064116e-06"
$values[0,3] = [double]"1.292544766329229e-05"
$values[0,4] = [double]"1.482538436903269e-06"
$values[0,5] = [double]"6.429087079595774e-06"
$values[0,6] = [double]"3.886815648002084e-06"
$values[0,7] = [double]"2.06332833840861e-06"
$values[0,8] = [double]"1.391168439113244e-06"
$values[0,9] = [double]"4.585446731653064e-06"
$values[0,10] = [double]"1.588222949067131e-05"
$values[0,11] = [double]"1.086943484551739e-06"
$values[0,12] = [double]"1.35740774567239e-05"
$values[0,13] = [double]"2.963883844131487e-06"
$values[0,14] = [double]"1.011899985314813e-05"
$values[0,15] = [double]"2.4522209969291e-06"
$values[0,16] = [double]"5.108278401166899e-06"
$values[0,17] = [double]"5.808940386486938e-06"
$values[0,18] = [double]"2.019680778175825e-06"
$values[0,19] = [double]"4.225659722578712e-06"
$values[0,20] = [double]"2.016910002566874e-06"
$values[0,21] = [double]"1.375267856928986e-06"
$values[0,22] = [double]"4.733199602924287e-06"
$values[0,23] = [double]"2.736995611485327e-06"
$values[0,24] = [double]"4.725129656435456e-06"
$values[0,25] = [double]"1.91368417290505e-06"
$values[0,26] = [double]"1.892778186629585e-06"
$values[0,27] = [double]"7.151719728426542e-06"
$values[0,28] = [double]"3.561140431429521e-07"
$values[0,29] = [double]"1.609477294550743e-06"
$values[0,30] = [double]"4.054021701449528e-06"
$values[0,31] = [double]"5.107952802063664e-06"
$values[0,32] = [double]"2.042376763711218e-06"
$values[0,33] = [double]"6.30601334705716e-07"
$values[0,34] = [double]"8.325200155923085e-07"
$values[0,35] = [double]"3.116919970125309e-06"
$values[0,36] = [double]"3.027505727004609e-06"
$values[0,37] = [double]"7.169658601924311e-06"
$values[0,38] = [double]"1.164288050858886e-06"
$values[0,39] = [double]"4.84859810967464e-06"
$values[0,40] = [double]"7.236487817863235e-07"
$values[0,41] = [double]"2.463952171183337e-07"
$values[0,42] = [double]"3.773708385779173e-07"
$values[0,43] = [double]"1.639219249227608e-06"
$values[0,44] = [double]"9.967619689632556e-07"
$values[0,45] = [double]"4.341798558016308e-06"
$values[0,46] = [double]"3.428516265557846e-06"
$values[0,47] = [double]"5.507138212124119e-08"
$values[0,48] = [double]"4.541603630059399e-06"
$values[0,49] = [double]"7.725677278358489e-06"
$values[0,50] = [double]"6.634530222981994e-07"
$values[0,51] = [double]"4.750691914523486e-06"
$values[0,52] = [double]"1.825248887143971e-06"
$values[0,53] = [double]"9.350912932859501e-07"
$values[0,54] = [double]"3.482692250145192e-07"
$values[0,55] = [double]"1.095726111088879e-05"
$values[0,56] = [double]"1.750473529682495e-06"
$values[0,57] = [double]"6.761524673493113e-06"
$values[0,58] = [double]"1.392564740854141e-06"
$values[0,59] = [double]"8.659665127197513e-07"
$values[0,60] = [double]"1.362502075608063e-06"
$values[0,61] = [double]"9.588906095814309e-07"
$values[0,62] = [double]"3.371091906956281e-06"
$values[0,63] = [double]"7.316690243897028e-06"
$values[0,64] = [double]"8.752997928240802e-06"
$values[0,65] = [double]"1.844082021307258e-06"
$values[0,66] = [double]"3.099869218203821e-06"
$values[0,67] = [double]"2.461940766806947e-06"
$values[0,68] = [double]"9.64128048508428e-07"
$values[0,69] = [double]"5.644833436235785e-07"
$values[0,70] = [double]"1.499257678005961e-07"
$values[0,71] = [double]"1.948806811924442e-06"
$values[0,72] = [double]"4.042138698423514e-06"
$values[0,73] = [double]"2.680922307263245e-06"
$values[0,74] = [double]"4.71225575893186e-07"
$values[0,75] = [double]"1.621771048121445e-06"
$values[0,76] = [double]"5.355992016120581e-06"
$values[0,77] = [double]"2.717900997595279e-07"
$values[0,78] = [double]"4.244430783728603e-06"
$values[0,79] = [double]"5.593310561380349e-07"
$values[0,80] = [double]"3.03219258057652e-06"
$values[0,81] = [double]"3.208181851732661e-06"
$values[0,82] = [double]"2.942622643331561e-07"
$values[0,83] = [double]"6.919680686223728e-07"
$values[0,84] = [double]"4.184933914075373e-06"
$values[0,85] = [double]"2.557798097768682e-06"
$values[0,86] = [double]"7.340219099205569e-07"
$values[0,87] = [double]"4.076101731698145e-07"
$values[0,88] = [double]"1.065903234120924e-06"
$values[0,89] = [double]"3.207886152267747e-07"
$values[0,90] = [double]"2.316873178642709e-06"
$values[0,91] = [double]"5.416384283307707e-07"
$values[0,92] = [double]"1.591911086507025e-06"
$values[0,93] = [double]"3.936450866603991e-06"
$values[0,94] = [double]"3.326214823573537e-07"
$values[0,95] = [double]"4.572078069031704e-06"
$values[0,96] = [double]"6.542798018926987e-06"
$values[0,97] = [double]"6.367437435983447e-06"
$values[0,98] = [double]"7.16824706614716e-06"
$values[0,99] = [double]"3.154159003315726e-06"
$values[0,100] = [double]"8.006355756151606e-07"
$values[0,101] = [double]"6.083071752982505e-07"
$values[0,102] = [double]"2.85050646198215e-06"
$values[0,103] = [double]"9.040613804245368e-07"
$values[0,104] = [double]"2.052836180155282e-06"
$values[0,105] = [double]"2.803642473736545e-07"
$values[0,106] = [double]"1.62117555646546e-06"
$values[0,107] = [double]"3.947154993966251e-07"
$values[0,108] = [double]"1.675614385021618e-06"
$values[0,109] = [double]"9.748262527864426e-06"
$values[0,110] = [double]"3.142667992506176e-06"
$values[0,111] = [double]"3.558825937943766e-06"
$values[0,112] = [double]"2.11526116800087e-06"
$values[0,113] = [double]"6.806977125961566e-06"
$values[0,114] = [double]"9.398782822245266e-06"
$values[0,115] = [double]"1.382921254844405e-06"
$values[0,116] = [double]"7.158431344578275e-06"
$values[0,117] = [double]"3.869455667881994e-06"
$values[0,118] = [double]"2.956134721898707e-06"
$values[0,119] = [double]"8.010476449271664e-06"
$values[0,120] = [double]"7.202733740996337e-06"
$values[0,121] = [double]"2.680332045201794e-06"
$values[0,122] = [double]"3.203614085123263e-08"
$values[0,123] = [double]"9.485962664257386e-07"
$values[0,124] = [double]"1.258311840501847e-06"
$values[0,125] = [double]"1.467611809857772e-06"
$values[0,126] = [double]"7.211657248262782e-06"
$values[0,127] = [double]"2.369937874391326e-06"
$values[0,128] = [double]"3.046949814233813e-06"
$values[0,129] = [double]"2.470600520609878e-06"
$values[0,130] = [double]"6.382758783729514e-07"
$values[0,131] = [double]"5.113337920192862e-06"
$values[0,132] = [double]"1.172140741800831e-06"
$values[0,133] = [double]"2.768649665085832e-06"
$values[0,134] = [double]"8.827050237414369e-07"
$values[0,135] = [double]"3.141248043903033e-06"
$values[0,136] = [double]"1.242074176843744e-06"
$values[0,137] = [double]"7.677658118154795e-07"
$values[0,138] = [double]"7.53190579416696e-06"
$values[0,139] = [double]"1.258206566490117e-06"
$values[0,140] = [double]"2.33203763855272e-06"
$values[0,141] = [double]"5.103941020934144e-07"
$values[0,142] = [double]"2.566808689152822e-06"
$values[0,143] = [double]"3.103066887888417e-07"
$values[0,144] = [double]"5.607028015219839e-06"
$values[0,145] = [double]"1.727313019728172e-06"
$values[0,146] = [double]"5.221010269451654e-06"
$values[0,147] = [double]"7.28020722817746e-07"
$values[0,148] = [double]"1.890717840069556e-06"
$values[0,149] = [double]"3.022542387043359e-06"
$values[0,150] = [double]"1.890856992758927e-06"
$values[0,151] = [double]"2.405248324066633e-06"
$values[0,152] = [double]"1.198633725607579e-07"
$values[0,153] = [double]"4.151949724473525e-06"
$values[0,154] = [double]"7.035915245978686e-07"
$values[0,155] = [double]"3.288046741545259e-07"
$values[0,156] = [double]"3.815157469944097e-06"
$values[0,157] = [double]"2.053349135167082e-06"
$values[0,158] = [double]"8.906915240913804e-07"
$values[0,159] = [double]"1.87288605957292e-07"
$values[0,160] = [double]"6.783445769542595e-07"
$values[0,161] = [double]"4.627851240002201e-07"
$values[0,162] = [double]"2.171806272599497e-06"
$values[0,163] = [double]"5.264115770842182e-06"
$values[0,164] = [double]"8.416179184678185e-07"
$values[0,165] = [double]"3.433837264310569e-06"
$values[0,166] = [double]"8.359752428077627e-07"
$values[0,167] = [double]"3.58756051355158e-06"
$values[0,168] = [double]"2.206753151767771e-06"
$values[0,169] = [double]"4.7910327793943e-07"
$values[0,170] = [double]"4.09509766541305e-07"
$values[0,171] = [double]"7.852395356167108e-06"
$values[0,172] = [double]"6.926562491571531e-06"
$values[0,173] = [double]"2.129381527993246e-06"
$values[0,174] = [double]"9.461054446546768e-07"
$values[0,175] = [double]"1.384051756758709e-06"
$values[0,176] = [double]"3.563448217391851e-06"
$values[0,177] = [double]"1.168573271570494e-05"
$values[0,178] = [double]"5.980455171084031e-06"
$values[0,179] = [double]"3.567659405234735e-06"
$values[0,180] = [double]"9.152076017926447e-06"
$values[0,181] = [double]"2.913708840424079e-06"
$values[0,182] = [double]"8.948064532887656e-06"
$values[0,183] = [double]"5.098399924463592e-08"
$values[0,184] = [double]"4.873204488831107e-06"
$values[0,185] = [double]"5.256324129732093e-06"
$values[0,186] = [double]"2.533623955969233e-06"
$values[0,187] = [double]"1.92257948583574e-06"
$values[0,188] = [double]"2.199493110310868e-06"

$ws.Range("A2:GG2").Value = $values
Write-Host "Done"